# Turn File Contents.xlsx - Property Tab tidy-up / Hull map editing re-enable
#
# This reproduces (at the cell-content / formatting level) the changes made to
# Sheet1 of "Turn File Contents.xlsx":
#   - "RaceData (not a RaceData object)" (D39) renamed to "PlayerRace"
#   - "BattlePlans" (K32) renamed to "BattlePlans (in PlayerData)"; duplicate L32 removed
#   - "PlayerRelations" (K34) renamed to "PlayerRelations (in PlayerData)"; duplicate L34 removed
#   - K39 ("RaceData" -> "PlayerRace") un-struck (no longer marked obsolete) and a
#     matching L39 entry added
#   - A new K47 "PlayerData" entry added (struck-through, matching its row siblings)
#   - The sheet's scroll/selection moved from D40 to I36

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 32: BattlePlans -----------------------------------------------
$ws.Range("K32").Value = "BattlePlans (in PlayerData)"
$ws.Range("L32").ClearContents()

# --- Row 34: PlayerRelations --------------------------------------------
$ws.Range("K34").Value = "PlayerRelations (in PlayerData)"
$ws.Range("L34").ClearContents()

# --- Row 39: RaceData -> PlayerRace --------------------------------------
$ws.Range("D39").Value = "PlayerRace"
$ws.Range("K39").Value = "PlayerRace"
$ws.Range("K39").Font.Strikethrough = $false
$ws.Range("L39").Value = "PlayerRace"

# --- Row 47: add PlayerData entry in column K ----------------------------
$ws.Range("K47").Value = "PlayerData"
$ws.Range("K47").Font.Strikethrough = $true

# --- Update the saved scroll position / selection ------------------------
$ws.Activate()
$ws.Range("I36").Select()
